$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated measurement values ---
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 10
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 9
$ws.Range("C4").Value = 1.4

# --- Column widths: column A and C widened, column B back to an
#     un-customized (auto) width. ColumnWidth is expressed in character
#     units and gets snapped to the same pixel grid real Excel uses, so
#     the inputs below are chosen to land on the saved widths 27 / 27.25. ---
$ws.Columns.Item(1).ColumnWidth = 26.3
$ws.Columns.Item(3).ColumnWidth = 26.5

# --- Selection moves to C4 ---
$ws.Range("C4").Select() | Out-Null
